# Update cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.304.67'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").Value = '1.872.16'
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("D4").Value = '''1.0000'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E5").Value = '  -1.57%  '
$ws.Range("D6").Value = '''241.64'
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").Value = '''0.9999'
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '''0.3108'
$ws.Range("E8").Value = '  +0.65%  '
$ws.Range("D9").Value = '''0.07695'
$ws.Range("E9").Value = '  -1.68%  '
$ws.Range("D10").Value = '''25.14'
$ws.Range("E10").Value = '  -0.26%  '
$ws.Range("D11").Value = '''0.08359'
$ws.Range("E11").Value = '  +1.32%  '
$ws.Range("D12").Value = '1.884.75'
$ws.Range("E12").Value = '  +1.40%  '
$ws.Range("D13").Value = '''5.222'
$ws.Range("E13").Value = '  -0.35%  '
$ws.Range("D14").Value = '''0.7113'
$ws.Range("E14").Value = '  -1.47%  '
$ws.Range("D15").Value = '''91.29'
$ws.Range("E15").Value = '  +0.58%  '
$ws.Range("D16").Value = '29.326.12'
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("D17").Value = '''0.000008257'
$ws.Range("E17").Value = '  +5.64%  '
$ws.Range("D18").Value = '''5.938'
$ws.Range("E18").Value = '  +1.36%  '
$ws.Range("D19").Value = '''242.21'
$ws.Range("E19").Value = '  -0.65%  '
$ws.Range("D20").Value = '2.133.41'
$ws.Range("E20").Value = '  +1.55%  '
$ws.Range("D21").Value = '''13.18'
$ws.Range("D22").Value = '''0.9995'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value = '''7.854'
$ws.Range("E23").Value = '  -1.68%  '
$ws.Range("E25").Value = '  +2.27%  '
$ws.Range("D26").Value = '''163.23'
$ws.Range("E26").Value = '  +0.77%  '
$ws.Range("D27").Value = '''9.012'
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("E28").Value = '  +1.43%  '
$ws.Range("E29").Value = '  +0.60%  '
$ws.Range("E30").Value = '  +0.19%  '
$ws.Range("D31").Value = '''4.342'
$ws.Range("E31").Value = '  +5.86%  '
$ws.Range("D32").Value = '''1.285'
$ws.Range("E32").Value = '  -4.60%  '
$ws.Range("D33").Value = '''0.05246'
$ws.Range("E33").Value = '  +0.95%  '
$ws.Range("D34").Value = '''1.928'
$ws.Range("E34").Value = '  -0.27%  '
$ws.Range("D35").Value = '''0.7542'
$ws.Range("E35").Value = '  +3.61%  '
$ws.Range("D36").Value = '''1.172'
$ws.Range("E36").Value = '  -1.27%  '
$ws.Range("D37").Value = '''2.679'
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").Value = '''0.01857'
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("D39").Value = '''2.718'
$ws.Range("E39").Value = '  +0.65%  '
$ws.Range("D40").Value = '1.152.66'
$ws.Range("E40").Value = '  -1.90%  '
$ws.Range("D41").Value = '''6.357'
$ws.Range("E41").Value = '  +4.27%  '
$ws.Range("D42").Value = '''73.10'
$ws.Range("E42").Value = '  +1.00%  '
$ws.Range("D43").Value = '''0.8886'
$ws.Range("E43").Value = '  -1.54%  '
$ws.Range("E44").Value = '  +2.88%  '
$ws.Range("D46").Value = '2.029.20'
$ws.Range("E46").Value = '  +1.50%  '
$ws.Range("D47").Value = '''0.5195'
$ws.Range("E47").Value = '  -1.69%  '
$ws.Range("D48").Value = '''1.793'
$ws.Range("E48").Value = '  +0.66%  '
$ws.Range("D49").Value = '''9.384'
$ws.Range("E49").Value = '  +0.95%  '
$ws.Range("E50").Value = '  +0.67%  '
$ws.Range("B51").Value = 'Frax'
$ws.Range("C51").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D51").Value = '''0.9978'
$ws.Range("E51").Value = '  +0.41%  '
